$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 47 for the Day-20 challenge entry.
# Shared-string insertion order must match: link, title, category, note.
$ws.Range("D47").Value = "https://leetcode.com/problems/short-encoding-of-words/"
$ws.Range("B47").Value = "Short Encoding of Words (June Day-20)"
$ws.Range("A47").Value = "47. String"
$ws.Range("C47").Value = "Use std::set"

# Row 46 currently holds the "Search Suggestions System (June Day-19)" entry;
# add its Category label in column A ("46. Trie")
$ws.Range("A46").Value = "46. Trie"

$ws.Hyperlinks.Add($ws.Range("D47"), "https://leetcode.com/problems/short-encoding-of-words/") | Out-Null

$ws.Range("B46").Copy() | Out-Null
$ws.Range("B47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("C46").Select()
